$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.279.17"
$ws.Range("E2").Value = "  +0.00%  "
$ws.Range("D3").Value = "1.667.05"
$ws.Range("E3").Value = "  +0.24%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'219.79"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("E6").Value = "  -0.48%  "
$ws.Range("E7").Value = "  -0.10%  "
$ws.Range("E8").Value = "  +0.32%  "
$ws.Range("D9").Value = "'0.06364"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").Value = "'20.91"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Value = "'0.07829"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.27%  "
$ws.Range("D12").Value = "'4.521"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.65%  "
$ws.Range("D13").Value = "1.664.76"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").Value = "1.895.91"
$ws.Range("E14").Value = "  +0.21%  "
$ws.Range("D15").Value = "'0.5600"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "0.0₅8111"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").Value = "'65.66"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "26.298.83"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").Value = "'4.721"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.28%  "
$ws.Range("D21").Value = "'199.64"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.70%  "
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").Value = "'6.046"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.12%  "
$ws.Range("D25").Value = "'146.12"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.1214"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.62%  "
$ws.Range("D27").Value = "'7.234"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("E28").Value = "  -0.21%  "
$ws.Range("D29").Value = "'1.529"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.07%  "
$ws.Range("D30").Value = "'0.05909"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("E31").Value = "  +0.21%  "
$ws.Range("D32").Value = "'3.514"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.07%  "
$ws.Range("D33").Value = "'3.321"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.58%  "
$ws.Range("D34").Value = "'1.597"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.07%  "
$ws.Range("D35").Value = "'0.9612"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.28%  "
$ws.Range("D36").Value = "'2.820"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "'2.429"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "'0.5790"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.32%  "
$ws.Range("D39").Value = "'0.01612"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").Value = "'5.956"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.22%  "
$ws.Range("D41").Value = "1.074.51"
$ws.Range("E41").Value = "  +2.59%  "
$ws.Range("D42").Value = "'0.8574"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.04%  "
$ws.Range("E43").Value = "  -0.07%  "
$ws.Range("D44").Value = "'102.75"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.50%  "
$ws.Range("D45").Value = "1.806.22"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").Value = "'58.44"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.99%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "0.0₈107"
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").Value = "'1.013"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("B49").Value = "Mantle"
$ws.Range("C49").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D49").Value = "'0.4413"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.93%  "
$ws.Range("D50").Value = "'8.073"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.05%  "
$ws.Range("D51").Value = "'0.05145"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.34%  "
